# Natmi following Dr Hou advice
# Update LR-pair edge-weight table: Ligand/Receptor-expressing cell counts
# changed from 1 to 3 for all rows, which cascades into the dependent
# average/total expression, specificity and edge-weight columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 26.29132066666667
$ws.Range("H2").Value = 78.87396200000001
$ws.Range("I2").Value = 0.1411782207947891
$ws.Range("J2").Value = 0.1411782207947891
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 17.16653
$ws.Range("N2").Value = 51.49959
$ws.Range("O2").Value = 0.0560345397128279
$ws.Range("P2").Value = 0.0560345397128279
$ws.Range("Q2").Value = 451.3307449639533
$ws.Range("R2").Value = 4061.97670467558
$ws.Range("S2").Value = 0.007910856619711993
$ws.Range("T2").Value = 0.007910856619711993
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 26.29132066666667
$ws.Range("H3").Value = 78.87396200000001
$ws.Range("I3").Value = 0.1411782207947891
$ws.Range("J3").Value = 0.1411782207947891
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 256.4443053333333
$ws.Range("N3").Value = 769.332916
$ws.Range("O3").Value = 0.8370788162388805
$ws.Range("P3").Value = 0.8370788162388805
$ws.Range("Q3").Value = 6742.259464659243
$ws.Range("R3").Value = 60680.33518193319
$ws.Range("S3").Value = 0.1181772979416133
$ws.Range("T3").Value = 0.1181772979416133
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 26.29132066666667
$ws.Range("H4").Value = 78.87396200000001
$ws.Range("I4").Value = 0.1411782207947891
$ws.Range("J4").Value = 0.1411782207947891
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 32.74538866666666
$ws.Range("N4").Value = 98.236166
$ws.Range("O4").Value = 0.1068866440482915
$ws.Range("P4").Value = 0.1068866440482915
$ws.Range("Q4").Value = 860.9195137899657
$ws.Range("R4").Value = 7748.275624109692
$ws.Range("S4").Value = 0.01509006623346373
$ws.Range("T4").Value = 0.01509006623346373
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 143.4723713333333
$ws.Range("H5").Value = 430.417114
$ws.Range("I5").Value = 0.7704129577533824
$ws.Range("J5").Value = 0.7704129577533824
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.16653
$ws.Range("N5").Value = 51.49959
$ws.Range("O5").Value = 0.0560345397128279
$ws.Range("P5").Value = 0.0560345397128279
$ws.Range("Q5").Value = 2462.922766664806
$ws.Range("R5").Value = 22166.30489998326
$ws.Range("S5").Value = 0.0431697354765091
$ws.Range("T5").Value = 0.0431697354765091
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 143.4723713333333
$ws.Range("H6").Value = 430.417114
$ws.Range("I6").Value = 0.7704129577533824
$ws.Range("J6").Value = 0.7704129577533824
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 256.4443053333333
$ws.Range("N6").Value = 769.332916
$ws.Range("O6").Value = 0.8370788162388805
$ws.Range("P6").Value = 0.8370788162388805
$ws.Range("Q6").Value = 36792.67260110271
$ws.Range("R6").Value = 331134.0534099244
$ws.Range("S6").Value = 0.6448963666912959
$ws.Range("T6").Value = 0.6448963666912959
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 143.4723713333333
$ws.Range("H7").Value = 430.417114
$ws.Range("I7").Value = 0.7704129577533824
$ws.Range("J7").Value = 0.7704129577533824
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 32.74538866666666
$ws.Range("N7").Value = 98.236166
$ws.Range("O7").Value = 0.1068866440482915
$ws.Range("P7").Value = 0.1068866440482915
$ws.Range("Q7").Value = 4698.058562238324
$ws.Range("R7").Value = 42282.52706014492
$ws.Range("S7").Value = 0.08234685558557724
$ws.Range("T7").Value = 0.08234685558557724
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 16.46418733333334
$ws.Range("H8").Value = 49.39256200000001
$ws.Range("I8").Value = 0.08840882145182853
$ws.Range("J8").Value = 0.08840882145182853
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 17.16653
$ws.Range("N8").Value = 51.49959
$ws.Range("O8").Value = 0.0560345397128279
$ws.Range("P8").Value = 0.0560345397128279
$ws.Range("Q8").Value = 282.6329657832867
$ws.Range("R8").Value = 2543.69669204958
$ws.Range("S8").Value = 0.004953947616606796
$ws.Range("T8").Value = 0.004953947616606796
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 16.46418733333334
$ws.Range("H9").Value = 49.39256200000001
$ws.Range("I9").Value = 0.08840882145182853
$ws.Range("J9").Value = 0.08840882145182853
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 256.4443053333333
$ws.Range("N9").Value = 769.332916
$ws.Range("O9").Value = 0.8370788162388805
$ws.Range("P9").Value = 0.8370788162388805
$ws.Range("Q9").Value = 4222.147083574532
$ws.Range("R9").Value = 37999.32375217079
$ws.Range("S9").Value = 0.07400515160597117
$ws.Range("T9").Value = 0.07400515160597117
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 16.46418733333334
$ws.Range("H10").Value = 49.39256200000001
$ws.Range("I10").Value = 0.08840882145182853
$ws.Range("J10").Value = 0.08840882145182853
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 32.74538866666666
$ws.Range("N10").Value = 98.236166
$ws.Range("O10").Value = 0.1068866440482915
$ws.Range("P10").Value = 0.1068866440482915
$ws.Range("Q10").Value = 539.1262133108103
$ws.Range("R10").Value = 4852.135919797292
$ws.Range("S10").Value = 0.009449722229250557
$ws.Range("T10").Value = 0.009449722229250557
